$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

# Raw and Clean data from SSA for September 11th-13th
$dates  = @("2020-09-11", "2020-09-12", "2020-09-13")
$values = @(
    @(658299, 750813, 87210, 70183, 24.73),
    @(663973, 759188, 87150, 70604, 24.67),
    @(668381, 765337, 82870, 70821, 24.59)
)

$startRow = 104
for ($i = 0; $i -lt $dates.Count; $i++) {
    $r = $startRow + $i

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dates[$i]

    for ($j = 0; $j -lt 5; $j++) {
        $cell = $ws.Cells.Item($r, 2 + $j)
        $cell.Value = $values[$i][$j]
        $cell.Font.Size = 12
        $cell.WrapText = $true
    }
}
